$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix phone number in C2 to be stored as a number instead of text
$ws.Range("C2").Value = 9112868820

# Add new member row 3
$ws.Range("A3").Value = "Arian"
$ws.Range("B3").Value = "Saeedkondori"
$ws.Range("C3").Value = "'09112868830"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "ariansk"
$ws.Range("E3").Value = "arian@gmail.com"
$ws.Range("F3").Value = "Arian@12"
$ws.Range("G3").Value = "Karaj"
$ws.Range("H3").Value = "2005/January/1"
$ws.Range("I3").Value = "city"
$ws.Range("J3").Value = "Karaj"
